$d = $word.ActiveDocument

function Find-ParagraphIndexStartingWith($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

function Merge-ParagraphRuns($paraIndex, $newText) {
    # Collapse all the runs of the paragraph (minus the trailing paragraph
    # mark) into a single run holding $newText, re-using Find/Replace so the
    # surviving run keeps its original run formatting (rPr).
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $rng = $d.Range($start, $end)
    $oldText = $rng.Text
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# The three entries below were each split across two runs (e.g. "02/03/2022 - "
# + "iterazione 3 / ..."); re-join each pair into a single run with the same
# visible text.
$idx0203 = Find-ParagraphIndexStartingWith "02/03/2022"
Merge-ParagraphRuns $idx0203 "02/03/2022 - iterazione 3 / stesura documentazione definitiva;"

$idx0303 = Find-ParagraphIndexStartingWith "03/03/2022"
Merge-ParagraphRuns $idx0303 "03/03/2022 - iterazione 3 / stesura documentazione definitiva;"

$idx0403 = Find-ParagraphIndexStartingWith "04/03/2022"
Merge-ParagraphRuns $idx0403 "04/03/2022 - iterazione 3 / stesura documentazione definitiva;"

# Add the new "project closure" entry right after the 04/03/2022 line, using
# the same paragraph formatting (hanging indent + it-IT language).
$enDash = [char]0x2013
$p0403 = $d.Paragraphs($idx0403)
$p0403.Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs($idx0403 + 1).Range.Text = "07/03/2022 $enDash chiusura progetto;"
